# Regenerate save_data to use K (strikeouts) instead of Strike# column values.
# The "K" column (column G) is recalculated/rewritten with updated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new "K" (strikeouts) value for column G.
$kValues = @{
    2  = 6
    3  = 5
    4  = 11
    5  = 3
    6  = 10
    7  = 8
    8  = 4
    9  = 9
    10 = 6
    11 = 10
    12 = 4
    13 = 7
    14 = 5
    15 = 5
    16 = 9
    17 = 10
    18 = 8
    19 = 6
    20 = 3
    21 = 3
    22 = 5
    23 = 12
    24 = 6
    25 = 6
    26 = 8
    27 = 5
    28 = 5
    29 = 3
    30 = 8
    31 = 5
    32 = 8
    33 = 7
    34 = 5
    35 = 6
    36 = 2
    37 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
